$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 0.2594103048008066
$ws.Range("J2").Value = 0.2594103048008066
$ws.Range("M2").Value = 227.11144
$ws.Range("N2").Value = 681.33432
$ws.Range("O2").Value = 0.8625743548356182
$ws.Range("P2").Value = 0.8625743548356182
$ws.Range("Q2").Value = 3.03852395576
$ws.Range("R2").Value = 27.34671560184
$ws.Range("S2").Value = 0.2237606763012668
$ws.Range("T2").Value = 0.2237606763012668

# Row 3
$ws.Range("I3").Value = 0.2594103048008066
$ws.Range("J3").Value = 0.2594103048008066
$ws.Range("O3").Value = 0.001598666154760757
$ws.Range("P3").Value = 0.001598666154760757
$ws.Range("S3").Value = 0.0004147104744812214
$ws.Range("T3").Value = 0.0004147104744812213

# Row 4
$ws.Range("I4").Value = 0.2594103048008066
$ws.Range("J4").Value = 0.2594103048008066
$ws.Range("M4").Value = 3.233093
$ws.Range("N4").Value = 9.699279000000001
$ws.Range("O4").Value = 0.01227935989749593
$ws.Range("P4").Value = 0.01227935989749593
$ws.Range("Q4").Value = 0.04325555124700001
$ws.Range("R4").Value = 0.389299961223
$ws.Range("S4").Value = 0.003185392493768221
$ws.Range("T4").Value = 0.00318539249376822

# Row 5
$ws.Range("I5").Value = 0.2594103048008066
$ws.Range("J5").Value = 0.2594103048008066
$ws.Range("M5").Value = 32.52945966666667
$ws.Range("N5").Value = 97.588379
$ws.Range("O5").Value = 0.1235476191121251
$ws.Range("P5").Value = 0.1235476191121251
$ws.Range("Q5").Value = 0.4352116408803334
$ws.Range("R5").Value = 3.916904767923
$ws.Range("S5").Value = 0.03204952553129035
$ws.Range("T5").Value = 0.03204952553129034

# Row 6
$ws.Range("G6").Value = 0.03819566666666666
$ws.Range("H6").Value = 0.114587
$ws.Range("I6").Value = 0.7405896951991934
$ws.Range("J6").Value = 0.7405896951991934
$ws.Range("M6").Value = 227.11144
$ws.Range("N6").Value = 681.33432
$ws.Range("O6").Value = 0.8625743548356182
$ws.Range("P6").Value = 0.8625743548356182
$ws.Range("Q6").Value = 8.674672858426666
$ws.Range("R6").Value = 78.07205572584
$ws.Range("S6").Value = 0.6388136785343513
$ws.Range("T6").Value = 0.6388136785343513

# Row 7
$ws.Range("G7").Value = 0.03819566666666666
$ws.Range("H7").Value = 0.114587
$ws.Range("I7").Value = 0.7405896951991934
$ws.Range("J7").Value = 0.7405896951991934
$ws.Range("O7").Value = 0.001598666154760757
$ws.Range("P7").Value = 0.001598666154760757
$ws.Range("Q7").Value = 0.01607734547711111
$ws.Range("R7").Value = 0.144696109294
$ws.Range("S7").Value = 0.001183955680279535
$ws.Range("T7").Value = 0.001183955680279535

# Row 8
$ws.Range("G8").Value = 0.03819566666666666
$ws.Range("H8").Value = 0.114587
$ws.Range("I8").Value = 0.7405896951991934
$ws.Range("J8").Value = 0.7405896951991934
$ws.Range("M8").Value = 3.233093
$ws.Range("N8").Value = 9.699279000000001
$ws.Range("O8").Value = 0.01227935989749593
$ws.Range("P8").Value = 0.01227935989749593
$ws.Range("Q8").Value = 0.1234901425303333
$ws.Range("R8").Value = 1.111411282773
$ws.Range("S8").Value = 0.009093967403727711
$ws.Range("T8").Value = 0.009093967403727709

# Row 9
$ws.Range("G9").Value = 0.03819566666666666
$ws.Range("H9").Value = 0.114587
$ws.Range("I9").Value = 0.7405896951991934
$ws.Range("J9").Value = 0.7405896951991934
$ws.Range("M9").Value = 32.52945966666667
$ws.Range("N9").Value = 97.588379
$ws.Range("O9").Value = 0.1235476191121251
$ws.Range("P9").Value = 0.1235476191121251
$ws.Range("Q9").Value = 1.242484398274778
$ws.Range("R9").Value = 11.182359584473
$ws.Range("S9").Value = 0.0914980935808348
$ws.Range("T9").Value = 0.0914980935808348
